{"js": "// 1. Clear the stray paragraph that contains only the single letter \"c\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.text === \"c\") {\n    p.clear();\n  }\n}\nawait context.sync();\n\n// 2. In the \"The file \"Examples\" ...\" paragraph, rename \"Examples\" -> \"Working Tutorial\".\nconst exampleHits = body.search(\"Examples\", { matchCase: true });\nexampleHits.load(\"text\");\nawait context.sync();\nif (exampleHits.items.length > 0) {\n  exampleHits.items[0].insertText(\"Working Tutorial\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. Extend \".. It illustrates most of the available operations.\" with an extra clause,\n//    right before the final period: \"... operations and allows you an environment to experiement.\"\nconst operationsHits = body.search(\"operations\", { matchCase: true });\noperationsHits.load(\"text\");\nawait context.sync();\nif (operationsHits.items.length > 0) {\n  const insertionPoint = operationsHits.items[0].getRange(\"End\");\n  insertionPoint.insertText(\n    \" and allows you an environment to experiement\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 4. Drop the leftover \"_GoBack\" bookmark (an internal last-edit marker Word had left behind).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Clear the stray paragraph that contains only the single letter \"c\".\n#    (Paragraph Range.Text includes the trailing paragraph mark, so trim it\n#    before comparing.)\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($paraText -eq \"c\") {\n        $p.Range.Text = \"\"\n    }\n}\n\n# 2. In the `The file \"Examples\" ...` paragraph, rename \"Examples\" -> \"Working Tutorial\".\n$find = $d.Content.Find\n$find.Text = \"Examples\"\n$find.Replacement.Text = \"Working Tutorial\"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# 3. Extend \"... It illustrates most of the available operations.\" with an extra\n#    clause, inserted right before the final period:\n#    \"... operations and allows you an environment to experiement.\"\n$find2 = $d.Content.Find\n$find2.Text = \"available operations.\"\n$find2.Execute() | Out-Null\n$matchRange = $find2.Parent\n$periodPos = $matchRange.End - 1\n$insertRange = $d.Range($periodPos, $periodPos)\n$insertRange.InsertBefore(\" and allows you an environment to experiement\")\n\n# 4. Remove the leftover \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
